# Auto-generated Excel COM-interop script
# Applies the cell-value updates described by the target diff
# (Sheets/Maduin_Profits.xlsx -> this workbook's ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR tabs).

$wb = $excel.ActiveWorkbook

# ----- ALC (38 cell updates) -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 275  # H11: 300 -> 275
$ws.Cells.Item(11, 9).Value = 275  # I11: 300 -> 275
$ws.Cells.Item(11, 11).Value = 275  # K11: 300 -> 275
$ws.Cells.Item(11, 13).Value = -135  # M11: -160 -> -135
$ws.Cells.Item(98, 8).Value = 456  # H98: 463 -> 456
$ws.Cells.Item(98, 9).Value = 456  # I98: 463 -> 456
$ws.Cells.Item(98, 11).Value = 456  # K98: 463 -> 456
$ws.Cells.Item(98, 13).Value = 1042  # M98: 1035 -> 1042
$ws.Cells.Item(122, 8).Value = 456  # H122: 463 -> 456
$ws.Cells.Item(122, 9).Value = 456  # I122: 463 -> 456
$ws.Cells.Item(122, 11).Value = 1368  # K122: 1389 -> 1368
$ws.Cells.Item(122, 13).Value = 1082  # M122: 1061 -> 1082
$ws.Cells.Item(125, 8).Value = 2185.9285  # H125: 2316.4614 -> 2185.9285
$ws.Cells.Item(125, 9).Value = 2088.875  # I125: 2202.625 -> 2088.875
$ws.Cells.Item(125, 10).Value = 2315.3333  # J125: 2498.6 -> 2315.3333
$ws.Cells.Item(125, 11).Value = 18799.875  # K125: 19823.625 -> 18799.875
$ws.Cells.Item(125, 12).Value = 20837.9997  # L125: 22487.4 -> 20837.9997
$ws.Cells.Item(125, 13).Value = -16339.875  # M125: -17363.625 -> -16339.875
$ws.Cells.Item(125, 14).Value = -25757.9997  # N125: -27407.4 -> -25757.9997
$ws.Cells.Item(131, 8).Value = 466.5  # H131: 447.5 -> 466.5
$ws.Cells.Item(131, 9).Value = 466.5  # I131: 447.5 -> 466.5
$ws.Cells.Item(131, 11).Value = 1399.5  # K131: 1342.5 -> 1399.5
$ws.Cells.Item(131, 13).Value = 3640.5  # M131: 3697.5 -> 3640.5
$ws.Cells.Item(135, 8).Value = 892  # H135: 900 -> 892
$ws.Cells.Item(135, 9).Value = 892  # I135: 900 -> 892
$ws.Cells.Item(135, 11).Value = 8028  # K135: 8100 -> 8028
$ws.Cells.Item(135, 13).Value = -5493  # M135: -5565 -> -5493
$ws.Cells.Item(137, 8).Value = 969.38464  # H137: 1009.8461 -> 969.38464
$ws.Cells.Item(137, 9).Value = 804.75  # I137: 837.44446 -> 804.75
$ws.Cells.Item(137, 10).Value = 1232.8  # J137: 1397.75 -> 1232.8
$ws.Cells.Item(137, 11).Value = 2414.25  # K137: 2512.33338 -> 2414.25
$ws.Cells.Item(137, 12).Value = 3698.4  # L137: 4193.25 -> 3698.4
$ws.Cells.Item(137, 13).Value = 135.75  # M137: 37.66661999999997 -> 135.75
$ws.Cells.Item(137, 14).Value = -8798.4  # N137: -9293.25 -> -8798.4
$ws.Cells.Item(141, 8).Value = 987.8  # H141: 991.2 -> 987.8
$ws.Cells.Item(141, 9).Value = 987.8  # I141: 991.2 -> 987.8
$ws.Cells.Item(141, 11).Value = 2963.4  # K141: 2973.6 -> 2963.4
$ws.Cells.Item(141, 13).Value = 2216.6  # M141: 2206.4 -> 2216.6

# ----- ARM (37 cell updates) -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 37342  # H43: 37341.5 -> 37342
$ws.Cells.Item(43, 10).Value = 0  # J43: 37341 -> 0
$ws.Cells.Item(43, 12).Value = 0  # L43: 37341 -> 0
$ws.Cells.Item(43, 14).ClearContents()  # N43: -37967 -> (cleared)
$ws.Cells.Item(61, 8).Value = 2154  # H61: 1989 -> 2154
$ws.Cells.Item(61, 9).Value = 2154  # I61: 1989 -> 2154
$ws.Cells.Item(61, 11).Value = 2154  # K61: 1989 -> 2154
$ws.Cells.Item(61, 13).Value = -1942  # M61: -1777 -> -1942
$ws.Cells.Item(63, 8).Value = 4650.4  # H63: 5450 -> 4650.4
$ws.Cells.Item(63, 9).Value = 3688  # I63: 4083.3333 -> 3688
$ws.Cells.Item(63, 10).Value = 8500  # J63: 7500 -> 8500
$ws.Cells.Item(63, 11).Value = 3688  # K63: 4083.3333 -> 3688
$ws.Cells.Item(63, 12).Value = 8500  # L63: 7500 -> 8500
$ws.Cells.Item(63, 13).Value = -3002  # M63: -3397.3333 -> -3002
$ws.Cells.Item(63, 14).Value = -9872  # N63: -8872 -> -9872
$ws.Cells.Item(66, 8).Value = 4650.4  # H66: 5450 -> 4650.4
$ws.Cells.Item(66, 9).Value = 3688  # I66: 4083.3333 -> 3688
$ws.Cells.Item(66, 10).Value = 8500  # J66: 7500 -> 8500
$ws.Cells.Item(66, 11).Value = 18440  # K66: 20416.6665 -> 18440
$ws.Cells.Item(66, 12).Value = 42500  # L66: 37500 -> 42500
$ws.Cells.Item(66, 13).Value = -15008  # M66: -16984.6665 -> -15008
$ws.Cells.Item(66, 14).Value = -49364  # N66: -44364 -> -49364
$ws.Cells.Item(97, 8).Value = 1335.7778  # H97: 1507.3334 -> 1335.7778
$ws.Cells.Item(97, 9).Value = 1072.2858  # I97: 1197.2632 -> 1072.2858
$ws.Cells.Item(97, 10).Value = 2258  # J97: 2685.6 -> 2258
$ws.Cells.Item(97, 11).Value = 1072.2858  # K97: 1197.2632 -> 1072.2858
$ws.Cells.Item(97, 12).Value = 2258  # L97: 2685.6 -> 2258
$ws.Cells.Item(97, 13).Value = -576.2858000000001  # M97: -701.2632000000001 -> -576.2858000000001
$ws.Cells.Item(97, 14).Value = -3250  # N97: -3677.6 -> -3250
$ws.Cells.Item(132, 8).Value = 3498.8  # H132: 3499 -> 3498.8
$ws.Cells.Item(132, 9).Value = 3498.8  # I132: 3499 -> 3498.8
$ws.Cells.Item(132, 11).Value = 10496.4  # K132: 10497 -> 10496.4
$ws.Cells.Item(132, 13).Value = -7966.400000000001  # M132: -7967 -> -7966.400000000001
$ws.Cells.Item(136, 8).Value = 2154  # H136: 1989 -> 2154
$ws.Cells.Item(136, 9).Value = 2154  # I136: 1989 -> 2154
$ws.Cells.Item(136, 11).Value = 6462  # K136: 5967 -> 6462
$ws.Cells.Item(136, 13).Value = -3912  # M136: -3417 -> -3912

# ----- BSM (23 cell updates) -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2289  # H86: 2646.8 -> 2289
$ws.Cells.Item(86, 9).Value = 3097.25  # I86: 3963 -> 3097.25
$ws.Cells.Item(86, 11).Value = 3097.25  # K86: 3963 -> 3097.25
$ws.Cells.Item(86, 13).Value = -1974.25  # M86: -2840 -> -1974.25
$ws.Cells.Item(89, 8).Value = 2289  # H89: 2646.8 -> 2289
$ws.Cells.Item(89, 9).Value = 3097.25  # I89: 3963 -> 3097.25
$ws.Cells.Item(89, 11).Value = 15486.25  # K89: 19815 -> 15486.25
$ws.Cells.Item(89, 13).Value = -9870.25  # M89: -14199 -> -9870.25
$ws.Cells.Item(94, 8).Value = 2035.9  # H94: 2319.875 -> 2035.9
$ws.Cells.Item(94, 9).Value = 2173.2222  # I94: 2587 -> 2173.2222
$ws.Cells.Item(94, 10).Value = 800  # J94: 450 -> 800
$ws.Cells.Item(94, 11).Value = 2173.2222  # K94: 2587 -> 2173.2222
$ws.Cells.Item(94, 12).Value = 800  # L94: 450 -> 800
$ws.Cells.Item(94, 13).Value = -1722.2222  # M94: -2136 -> -1722.2222
$ws.Cells.Item(94, 14).Value = -1702  # N94: -1352 -> -1702
$ws.Cells.Item(99, 8).Value = 3544.3333  # H99: 4285.2856 -> 3544.3333
$ws.Cells.Item(99, 9).Value = 3487.5  # I99: 4333 -> 3487.5
$ws.Cells.Item(99, 11).Value = 3487.5  # K99: 4333 -> 3487.5
$ws.Cells.Item(99, 13).Value = -1989.5  # M99: -2835 -> -1989.5
$ws.Cells.Item(105, 8).Value = 3033  # H105: 3634.3333 -> 3033
$ws.Cells.Item(105, 9).Value = 2474.625  # I105: 3282.9092 -> 2474.625
$ws.Cells.Item(105, 11).Value = 2474.625  # K105: 3282.9092 -> 2474.625
$ws.Cells.Item(105, 13).Value = -727.625  # M105: -1535.9092 -> -727.625

# ----- CRP (25 cell updates) -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 3816.6667  # H22: 348.66666 -> 3816.6667
$ws.Cells.Item(22, 9).Value = 725  # I22: 348.66666 -> 725
$ws.Cells.Item(22, 10).Value = 10000  # J22: 0 -> 10000
$ws.Cells.Item(22, 11).Value = 725  # K22: 348.66666 -> 725
$ws.Cells.Item(22, 12).Value = 10000  # L22: 0 -> 10000
$ws.Cells.Item(22, 13).Value = -375  # M22: 1.333340000000021 -> -375
$ws.Cells.Item(22, 14).Value = -10700  # N22: (empty) -> -10700
$ws.Cells.Item(99, 8).Value = 3713.6667  # H99: 3260.1365 -> 3713.6667
$ws.Cells.Item(99, 9).Value = 3344.2354  # I99: 3016.45 -> 3344.2354
$ws.Cells.Item(99, 10).Value = 9994  # J99: 5697 -> 9994
$ws.Cells.Item(99, 11).Value = 3344.2354  # K99: 3016.45 -> 3344.2354
$ws.Cells.Item(99, 12).Value = 9994  # L99: 5697 -> 9994
$ws.Cells.Item(99, 13).Value = -1846.2354  # M99: -1518.45 -> -1846.2354
$ws.Cells.Item(99, 14).Value = -12990  # N99: -8693 -> -12990
$ws.Cells.Item(126, 8).Value = 3713.6667  # H126: 3260.1365 -> 3713.6667
$ws.Cells.Item(126, 9).Value = 3344.2354  # I126: 3016.45 -> 3344.2354
$ws.Cells.Item(126, 10).Value = 9994  # J126: 5697 -> 9994
$ws.Cells.Item(126, 11).Value = 10032.7062  # K126: 9049.349999999999 -> 10032.7062
$ws.Cells.Item(126, 12).Value = 29982  # L126: 17091 -> 29982
$ws.Cells.Item(126, 13).Value = -7562.706200000001  # M126: -6579.349999999999 -> -7562.706200000001
$ws.Cells.Item(126, 14).Value = -34922  # N126: -22031 -> -34922
$ws.Cells.Item(132, 8).Value = 3252.1428  # H132: 4995 -> 3252.1428
$ws.Cells.Item(132, 9).Value = 3252.1428  # I132: 4995 -> 3252.1428
$ws.Cells.Item(132, 11).Value = 9756.428400000001  # K132: 14985 -> 9756.428400000001
$ws.Cells.Item(132, 13).Value = -7226.428400000001  # M132: -12455 -> -7226.428400000001

# ----- CUL (26 cell updates) -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 273.66666  # H7: 228.375 -> 273.66666
$ws.Cells.Item(7, 9).Value = 86.666664  # I7: 86.5 -> 86.666664
$ws.Cells.Item(7, 10).Value = 460.66666  # J7: 370.25 -> 460.66666
$ws.Cells.Item(7, 11).Value = 259.999992  # K7: 259.5 -> 259.999992
$ws.Cells.Item(7, 12).Value = 1381.99998  # L7: 1110.75 -> 1381.99998
$ws.Cells.Item(7, 13).Value = -147.999992  # M7: -147.5 -> -147.999992
$ws.Cells.Item(7, 14).Value = -1605.99998  # N7: -1334.75 -> -1605.99998
$ws.Cells.Item(12, 8).Value = 92.333336  # H12: 97.59999999999999 -> 92.333336
$ws.Cells.Item(12, 10).Value = 92.333336  # J12: 97.59999999999999 -> 92.333336
$ws.Cells.Item(12, 12).Value = 277.000008  # L12: 292.8 -> 277.000008
$ws.Cells.Item(12, 14).Value = -623.000008  # N12: -638.8 -> -623.000008
$ws.Cells.Item(108, 8).Value = 373.66666  # H108: 323 -> 373.66666
$ws.Cells.Item(108, 9).Value = 373.66666  # I108: 323 -> 373.66666
$ws.Cells.Item(108, 11).Value = 1120.99998  # K108: 969 -> 1120.99998
$ws.Cells.Item(108, 13).Value = 1759.00002  # M108: 1911 -> 1759.00002
$ws.Cells.Item(129, 8).Value = 1578.091  # H129: 1741.8889 -> 1578.091
$ws.Cells.Item(129, 10).Value = 1845.25  # J129: 2180 -> 1845.25
$ws.Cells.Item(129, 12).Value = 5535.75  # L129: 6540 -> 5535.75
$ws.Cells.Item(129, 14).Value = -15535.75  # N129: -16540 -> -15535.75
$ws.Cells.Item(131, 8).Value = 971.1111  # H131: 925.5263 -> 971.1111
$ws.Cells.Item(131, 9).Value = 820  # I131: 700.8333 -> 820
$ws.Cells.Item(131, 10).Value = 1092  # J131: 1029.2307 -> 1092
$ws.Cells.Item(131, 11).Value = 2460  # K131: 2102.4999 -> 2460
$ws.Cells.Item(131, 12).Value = 3276  # L131: 3087.6921 -> 3276
$ws.Cells.Item(131, 13).Value = 2580  # M131: 2937.5001 -> 2580
$ws.Cells.Item(131, 14).Value = -13356  # N131: -13167.6921 -> -13356

# ----- GSM (30 cell updates) -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4499.1665  # H70: 4599.4 -> 4499.1665
$ws.Cells.Item(70, 9).Value = 4499.1665  # I70: 4599.4 -> 4499.1665
$ws.Cells.Item(70, 11).Value = 4499.1665  # K70: 4599.4 -> 4499.1665
$ws.Cells.Item(70, 13).Value = -4229.1665  # M70: -4329.4 -> -4229.1665
$ws.Cells.Item(73, 8).Value = 4499.1665  # H73: 4599.4 -> 4499.1665
$ws.Cells.Item(73, 9).Value = 4499.1665  # I73: 4599.4 -> 4499.1665
$ws.Cells.Item(73, 11).Value = 4499.1665  # K73: 4599.4 -> 4499.1665
$ws.Cells.Item(73, 13).Value = -3563.1665  # M73: -3663.4 -> -3563.1665
$ws.Cells.Item(96, 8).Value = 24999.5  # H96: 24998 -> 24999.5
$ws.Cells.Item(96, 10).Value = 24999.5  # J96: 24998 -> 24999.5
$ws.Cells.Item(96, 12).Value = 24999.5  # L96: 24998 -> 24999.5
$ws.Cells.Item(96, 14).Value = -30491.5  # N96: -30490 -> -30491.5
$ws.Cells.Item(122, 8).Value = 1300  # H122: 0 -> 1300
$ws.Cells.Item(122, 9).Value = 1300  # I122: 0 -> 1300
$ws.Cells.Item(122, 11).Value = 3900  # K122: 0 -> 3900
$ws.Cells.Item(122, 13).Value = -1450  # M122: (empty) -> -1450
$ws.Cells.Item(126, 8).Value = 4780  # H126: 4980 -> 4780
$ws.Cells.Item(126, 9).Value = 4000  # I126: 4333.3335 -> 4000
$ws.Cells.Item(126, 10).Value = 5300  # J126: 5950 -> 5300
$ws.Cells.Item(126, 11).Value = 12000  # K126: 13000.0005 -> 12000
$ws.Cells.Item(126, 12).Value = 15900  # L126: 17850 -> 15900
$ws.Cells.Item(126, 13).Value = -9530  # M126: -10530.0005 -> -9530
$ws.Cells.Item(126, 14).Value = -20840  # N126: -22790 -> -20840
$ws.Cells.Item(132, 8).Value = 3112.2727  # H132: 3515.375 -> 3112.2727
$ws.Cells.Item(132, 9).Value = 3349.25  # I132: 3849.3333 -> 3349.25
$ws.Cells.Item(132, 10).Value = 2480.3333  # J132: 2513.5 -> 2480.3333
$ws.Cells.Item(132, 11).Value = 10047.75  # K132: 11547.9999 -> 10047.75
$ws.Cells.Item(132, 12).Value = 7440.999899999999  # L132: 7540.5 -> 7440.999899999999
$ws.Cells.Item(132, 13).Value = -7517.75  # M132: -9017.999899999999 -> -7517.75
$ws.Cells.Item(132, 14).Value = -12500.9999  # N132: -12600.5 -> -12500.9999

# ----- LTW (20 cell updates) -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2267.8333  # H7: 2581.4 -> 2267.8333
$ws.Cells.Item(7, 9).Value = 2267.8333  # I7: 2581.4 -> 2267.8333
$ws.Cells.Item(7, 11).Value = 2267.8333  # K7: 2581.4 -> 2267.8333
$ws.Cells.Item(7, 13).Value = -2155.8333  # M7: -2469.4 -> -2155.8333
$ws.Cells.Item(40, 8).Value = 1896.8  # H40: 1780.5714 -> 1896.8
$ws.Cells.Item(40, 9).Value = 1746  # I40: 1660.6666 -> 1746
$ws.Cells.Item(40, 11).Value = 1746  # K40: 1660.6666 -> 1746
$ws.Cells.Item(40, 13).Value = -1610  # M40: -1524.6666 -> -1610
$ws.Cells.Item(46, 8).Value = 4430.769  # H46: 4400.0435 -> 4430.769
$ws.Cells.Item(46, 10).Value = 5120  # J46: 5022.3335 -> 5120
$ws.Cells.Item(46, 12).Value = 5120  # L46: 5022.3335 -> 5120
$ws.Cells.Item(46, 14).Value = -5496  # N46: -5398.3335 -> -5496
$ws.Cells.Item(126, 8).Value = 2267.8333  # H126: 2581.4 -> 2267.8333
$ws.Cells.Item(126, 9).Value = 2267.8333  # I126: 2581.4 -> 2267.8333
$ws.Cells.Item(126, 11).Value = 6803.499899999999  # K126: 7744.200000000001 -> 6803.499899999999
$ws.Cells.Item(126, 13).Value = -4333.499899999999  # M126: -5274.200000000001 -> -4333.499899999999
$ws.Cells.Item(132, 8).Value = 5375.8  # H132: 4222.25 -> 5375.8
$ws.Cells.Item(132, 9).Value = 4219.75  # I132: 3396.8572 -> 4219.75
$ws.Cells.Item(132, 11).Value = 12659.25  # K132: 10190.5716 -> 12659.25
$ws.Cells.Item(132, 13).Value = -10129.25  # M132: -7660.571599999999 -> -10129.25

# ----- WVR (18 cell updates) -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1486.9  # H126: 1484.5 -> 1486.9
$ws.Cells.Item(126, 9).Value = 1486.9  # I126: 1540.9 -> 1486.9
$ws.Cells.Item(126, 10).Value = 0  # J126: 1202.5 -> 0
$ws.Cells.Item(126, 11).Value = 4460.700000000001  # K126: 4622.700000000001 -> 4460.700000000001
$ws.Cells.Item(126, 12).Value = 0  # L126: 3607.5 -> 0
$ws.Cells.Item(126, 13).Value = -1990.700000000001  # M126: -2152.700000000001 -> -1990.700000000001
$ws.Cells.Item(126, 14).ClearContents()  # N126: -8547.5 -> (cleared)
$ws.Cells.Item(132, 8).Value = 1686  # H132: 1545.8 -> 1686
$ws.Cells.Item(132, 9).Value = 1640.4286  # I132: 1494.7778 -> 1640.4286
$ws.Cells.Item(132, 11).Value = 4921.2858  # K132: 4484.3334 -> 4921.2858
$ws.Cells.Item(132, 13).Value = -2391.2858  # M132: -1954.3334 -> -2391.2858
$ws.Cells.Item(136, 8).Value = 34081.066  # H136: 42494.5 -> 34081.066
$ws.Cells.Item(136, 9).Value = 39247.77  # I136: 46312.184 -> 39247.77
$ws.Cells.Item(136, 10).Value = 497.5  # J136: 500 -> 497.5
$ws.Cells.Item(136, 11).Value = 117743.31  # K136: 138936.552 -> 117743.31
$ws.Cells.Item(136, 12).Value = 1492.5  # L136: 1500 -> 1492.5
$ws.Cells.Item(136, 13).Value = -115193.31  # M136: -136386.552 -> -115193.31
$ws.Cells.Item(136, 14).Value = -6592.5  # N136: -6600 -> -6592.5

Write-Host "Applied 217 cell updates across 8 sheets"
